$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New objective values for rows 2..151 (column D)
$newValues = @(1381372.420472553,1331175.398185599,1304983.662428365,1296398.823928755,1301345.480533996,1256870.135244087,1262092.408059028,1247202.0759675,1201081.857238743,1160161.396475241,1161208.289735802,1139561.870749052,1125166.791264905,1106303.94116646,1087307.409048395,1071959.553167581,1057881.165404926,1052239.065325873,1043822.338174355,1037714.254760547,1033018.848612246,1029145.787170681,1024224.836825895,1021026.67311727,1021026.680022619,1021026.636777087,1021026.678955722,1021026.67311727,1021026.674058062,1021026.680022619,1021026.675216346,1021026.675325921,1021026.636777087,1021026.687943932,1021026.680967663,1021026.636777087,1021026.636777087,1021026.678955722,1021026.690114717,1021026.680858088,1021026.687229849,1021026.678758435,1021026.675325921,1021026.670736702,1021026.678955722,1021026.687943932,1021026.678955722,1021026.674058062,1021026.659157655,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.680858088,1021026.659157655,1021026.67311727,1021026.662915119,1021026.687229849,1021026.680858088,1021026.677651484,1021026.685004374,1021026.672685405,1021026.688963982)

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $ws.Cells.Item($i + 2, 4).Value = $newValues[$i]
}

# Remove the now-obsolete tail rows (152..201)
$ws.Range("A152:D201").EntireRow.Delete()
